$wb = $excel.ActiveWorkbook

# --- Update workbook-level absolute path metadata (author's environment path) ---
# (Not directly exposed via COM object model in a meaningful way; skipped.)

# --- Get sheets ---
$wsInner = $wb.Worksheets.Item("INNERMODEL")
$wsInner.Activate()

# Clear the old matrix content
$wsInner.Cells.Clear()

# Build the new FROM/TO connection list for the inner model
$pairs = @(
    @("FROM", "TO"),
    @("IMAGE", "EXPECTATION"),
    @("IMAGE", "SATISFACTION"),
    @("IMAGE", "LOYALTY"),
    @("EXPECTATION", "VALUE"),
    @("EXPECTATION", "QUALITY"),
    @("EXPECTATION", "SATISFACTION"),
    @("QUALITY", "VALUE"),
    @("QUALITY", "SATISFACTION"),
    @("VALUE", "SATISFACTION"),
    @("SATISFACTION", "LOYALTY")
)

for ($i = 0; $i -lt $pairs.Count; $i++) {
    $r = $i + 1
    $wsInner.Cells.Item($r, 1).Value = $pairs[$i][0]
    $wsInner.Cells.Item($r, 2).Value = $pairs[$i][1]
}

$wsInner.Range("C1:H10").Select()
